# Generate Report for Handback
#
# Adds "Latest Target File" (F) and "Latest Handback File" (G) hyperlinked
# entries to the zh-cn and de-de report rows, mirroring the handoff file
# links (source .md -> F, handoff .xlf -> G) now that the handback is in,
# updates the "Latest Handback DateTime" (H) timestamps, and flips the
# Status text from "Ready for handoff" to "Handed back: in sync with en-US"
# everywhere it appears (Overview sheet + per-language sheets).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---- Overview sheet: Status columns (B/C) for both language rows ----
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")

# Status column
$zh.Range("C2").Value = $statusText
$zh.Range("C3").Value = $statusText

# Row 2 (03b3fade...) - Latest Target File / Latest Handback File
$zh.Hyperlinks.Add($zh.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3cddaf00a6de69c537f8b66231940d9dac114da8/e2e/03b3fade-3891-462c-8010-5b9d4da9e099.md", "", "", "03b3fade-3891-462c-8010-5b9d4da9e099.md") | Out-Null
$zh.Range("F2").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43c7e47138e56509c029d94ad5a709ace8bc3468/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/03b3fade-3891-462c-8010-5b9d4da9e099.e68a72ea866663c2125a7c4278cc9dcb19f293af.zh-cn.xlf", "", "", "03b3fade-3891-462c-8010-5b9d4da9e099.e68a72ea866663c2125a7c4278cc9dcb19f293af.zh-cn.xlf") | Out-Null
$zh.Range("G2").Style = "HyperLink"
$zh.Range("H2").Value = "2016-03-23 03:04:28"

# Row 3 (eb823e95...) - Latest Target File / Latest Handback File
$zh.Hyperlinks.Add($zh.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3cddaf00a6de69c537f8b66231940d9dac114da8/e2e/eb823e95-f47d-40f6-8c6e-3557cda239d0.md", "", "", "eb823e95-f47d-40f6-8c6e-3557cda239d0.md") | Out-Null
$zh.Range("F3").Style = "HyperLink"
$zh.Hyperlinks.Add($zh.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/43c7e47138e56509c029d94ad5a709ace8bc3468/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/eb823e95-f47d-40f6-8c6e-3557cda239d0.49503e8115871fdf9dab20e72e60f7e06adbebc1.zh-cn.xlf", "", "", "eb823e95-f47d-40f6-8c6e-3557cda239d0.49503e8115871fdf9dab20e72e60f7e06adbebc1.zh-cn.xlf") | Out-Null
$zh.Range("G3").Style = "HyperLink"
$zh.Range("H3").Value = "2016-03-23 03:04:28"

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")

# Status column
$de.Range("C2").Value = $statusText
$de.Range("C3").Value = $statusText

# Row 2 (03b3fade...) - Latest Target File / Latest Handback File
$de.Hyperlinks.Add($de.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/3cddaf00a6de69c537f8b66231940d9dac114da8/e2e/03b3fade-3891-462c-8010-5b9d4da9e099.md", "", "", "03b3fade-3891-462c-8010-5b9d4da9e099.md") | Out-Null
$de.Range("F2").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b1fe52b96f06b37476899ee57ed95b50e0e9161/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/03b3fade-3891-462c-8010-5b9d4da9e099.e68a72ea866663c2125a7c4278cc9dcb19f293af.de-de.xlf", "", "", "03b3fade-3891-462c-8010-5b9d4da9e099.e68a72ea866663c2125a7c4278cc9dcb19f293af.de-de.xlf") | Out-Null
$de.Range("G2").Style = "HyperLink"
$de.Range("H2").Value = "2016-03-23 03:04:35"

# Row 3 (eb823e95...) - Latest Target File / Latest Handback File
$de.Hyperlinks.Add($de.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/3cddaf00a6de69c537f8b66231940d9dac114da8/e2e/eb823e95-f47d-40f6-8c6e-3557cda239d0.md", "", "", "eb823e95-f47d-40f6-8c6e-3557cda239d0.md") | Out-Null
$de.Range("F3").Style = "HyperLink"
$de.Hyperlinks.Add($de.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3b1fe52b96f06b37476899ee57ed95b50e0e9161/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/eb823e95-f47d-40f6-8c6e-3557cda239d0.49503e8115871fdf9dab20e72e60f7e06adbebc1.de-de.xlf", "", "", "eb823e95-f47d-40f6-8c6e-3557cda239d0.49503e8115871fdf9dab20e72e60f7e06adbebc1.de-de.xlf") | Out-Null
$de.Range("G3").Style = "HyperLink"
$de.Range("H3").Value = "2016-03-23 03:04:35"
